# New weekly price record for "Betarraga" at Macroferia Regional de Talca.
# A new row is inserted at row 94 (pushing the existing rows 94-187 down to
# 95-188), and the new row 94 is populated with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(94).Insert()

$ws.Range("A94").Value = 5
$ws.Range("B94").Value = "Macroferia Regional de Talca"
$ws.Range("C94").Value = "Maule"
$ws.Range("D94").Value = (Get-Date -Year 2021 -Month 10 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D94").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E94").Value = 7
$ws.Range("F94").Value = 100114014
$ws.Range("G94").Value = "Betarraga"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 700
$ws.Range("L94").Value = 700
$ws.Range("M94").Value = 700
$ws.Range("N94").Value = "`$/paquete 5 unidades"
$ws.Range("O94").Value = "Región del Maule"
$ws.Range("P94").Value = 140
$ws.Range("Q94").Value = 5
$ws.Range("R94").Value = "Hortaliza"
